# Auto-generated Excel COM-interop edit script
# Applies numeric cell-value updates per the target diff, sheet by sheet, row by row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1785.9048
$ws.Range("J17").Value = 1785.9048
$ws.Range("L17").Value = 5357.7144
$ws.Range("N17").Value = -5693.7144

$ws.Range("H18").Value = 468.75
$ws.Range("I18").Value = 468.75
$ws.Range("K18").Value = 468.75
$ws.Range("M18").Value = -184.75

$ws.Range("H34").Value = 1500
$ws.Range("I34").Value = 1500
$ws.Range("K34").Value = 1500
$ws.Range("M34").Value = -1297

$ws.Range("H36").Value = 1500
$ws.Range("I36").Value = 1500
$ws.Range("K36").Value = 1500
$ws.Range("M36").Value = -785

$ws.Range("H115").Value = 2262
$ws.Range("I115").Value = 2262
$ws.Range("K115").Value = 6786
$ws.Range("M115").Value = -5219

$ws.Range("H118").Value = 913.3333
$ws.Range("I118").Value = 913.3333
$ws.Range("K118").Value = 2739.9999
$ws.Range("M118").Value = -1082.9999

$ws.Range("H127").Value = 775
$ws.Range("J127").Value = 900
$ws.Range("L127").Value = 2700
$ws.Range("N127").Value = -12620

$ws.Range("H129").Value = 5700
$ws.Range("I129").Value = 5700
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 17100
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -12100

$ws.Range("H141").Value = 6999.5
$ws.Range("I141").Value = 6999.5
$ws.Range("K141").Value = 20998.5
$ws.Range("M141").Value = -15818.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4459.4053
$ws.Range("I32").Value = 4569.3887
$ws.Range("K32").Value = 4569.3887
$ws.Range("M32").Value = -4282.3887

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H80").Value = 121.1
$ws.Range("I80").Value = 103.2
$ws.Range("J80").Value = 139
$ws.Range("K80").Value = 103.2
$ws.Range("L80").Value = 139
$ws.Range("M80").Value = 894.8
$ws.Range("N80").Value = -2135

$ws.Range("H83").Value = 121.1
$ws.Range("I83").Value = 103.2
$ws.Range("J83").Value = 139
$ws.Range("K83").Value = 516
$ws.Range("L83").Value = 695
$ws.Range("M83").Value = 4476
$ws.Range("N83").Value = -10679

$ws.Range("H134").Value = 6480.1816
$ws.Range("I134").Value = 2612.1428
$ws.Range("J134").Value = 13249.25
$ws.Range("K134").Value = 7836.428400000001
$ws.Range("L134").Value = 39747.75
$ws.Range("M134").Value = -5301.428400000001
$ws.Range("N134").Value = -44817.75

$ws.Range("H141").Value = 159988
$ws.Range("J141").Value = 159987.5
$ws.Range("L141").Value = 159987.5
$ws.Range("N141").Value = -170347.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 375.78946
$ws.Range("I22").Value = 384.70587
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 384.70587
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -34.70587
$ws.Range("N22").Value = -1000

$ws.Range("H31").Value = 1485
$ws.Range("I31").Value = 1580
$ws.Range("K31").Value = 1580
$ws.Range("M31").Value = -1285

$ws.Range("H34").Value = 1485
$ws.Range("I34").Value = 1580
$ws.Range("K34").Value = 1580
$ws.Range("M34").Value = -1378

$ws.Range("H58").Value = 3185.276
$ws.Range("I58").Value = 3240.9583
$ws.Range("K58").Value = 3240.9583
$ws.Range("M58").Value = -3037.9583

$ws.Range("H99").Value = 2667
$ws.Range("I99").Value = 2667
$ws.Range("K99").Value = 2667
$ws.Range("M99").Value = -1169

$ws.Range("H126").Value = 2667
$ws.Range("I126").Value = 2667
$ws.Range("K126").Value = 8001
$ws.Range("M126").Value = -5531

$ws.Range("H136").Value = 3185.276
$ws.Range("I136").Value = 3240.9583
$ws.Range("K136").Value = 9722.874899999999
$ws.Range("M136").Value = -7172.874899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2667
$ws.Range("I5").Value = 2667
$ws.Range("K5").Value = 8001
$ws.Range("M5").Value = -7889

$ws.Range("H131").Value = 2149.6924
$ws.Range("J131").Value = 2157.9375
$ws.Range("L131").Value = 6473.8125
$ws.Range("N131").Value = -16553.8125

$ws.Range("H135").Value = 2667
$ws.Range("I135").Value = 2667
$ws.Range("K135").Value = 24003
$ws.Range("M135").Value = -21468

$ws.Range("H137").Value = 1399.5
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2469.9092
$ws.Range("I102").Value = 2469.9092
$ws.Range("K102").Value = 2469.9092
$ws.Range("M102").Value = -847.9092000000001

$ws.Range("H122").Value = 4358.857
$ws.Range("I122").Value = 3802.4
$ws.Range("J122").Value = 5750
$ws.Range("K122").Value = 11407.2
$ws.Range("L122").Value = 17250
$ws.Range("M122").Value = -8957.200000000001
$ws.Range("N122").Value = -22150

$ws.Range("H132").Value = 3685.4285
$ws.Range("I132").Value = 2750
$ws.Range("K132").Value = 8250
$ws.Range("M132").Value = -5720

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10429.733
$ws.Range("I22").Value = 13578
$ws.Range("J22").Value = 7675
$ws.Range("K22").Value = 13578
$ws.Range("L22").Value = 7675
$ws.Range("M22").Value = -13283
$ws.Range("N22").Value = -8265

$ws.Range("H27").Value = 10429.733
$ws.Range("I27").Value = 13578
$ws.Range("J27").Value = 7675
$ws.Range("K27").Value = 13578
$ws.Range("L27").Value = 7675
$ws.Range("M27").Value = -13471
$ws.Range("N27").Value = -7889

$ws.Range("H55").Value = 2346.6667
$ws.Range("I55").Value = 3082.5
$ws.Range("J55").Value = 875
$ws.Range("K55").Value = 3082.5
$ws.Range("L55").Value = 875
$ws.Range("M55").Value = -2909.5
$ws.Range("N55").Value = -1221

$ws.Range("H93").Value = 599.6667
$ws.Range("I93").Value = 599.6667
$ws.Range("K93").Value = 599.6667
$ws.Range("M93").Value = 648.3333

$ws.Range("H122").Value = 2628.4285
$ws.Range("I122").Value = 2779.8
$ws.Range("K122").Value = 8339.400000000001
$ws.Range("M122").Value = -5889.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3166.6667
$ws.Range("J126").Value = 2875
$ws.Range("L126").Value = 8625
$ws.Range("N126").Value = -13565

$ws.Range("H132").Value = 3582.8333
$ws.Range("I132").Value = 2499.5
$ws.Range("K132").Value = 7498.5
$ws.Range("M132").Value = -4968.5
